$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = '[-, -, -, -]'
$ws.Range("C3").Value = '-'
$ws.Range("D3").Value = '-'
$ws.Range("E3").Value = '[-, ''MCT-3A-Processos de Usinagem 2'', -, -]'
$ws.Range("F3").Value = '[-, ''MEC-3B-Retífica'', -, -]'
$ws.Range("C4").Value = '-'
$ws.Range("D4").Value = '-'
$ws.Range("E4").Value = '[-, ''MCT-3A-Processos de Usinagem 2'', -, -]'
$ws.Range("F4").Value = '[-, ''MEC-3B-Retífica'', -, -]'
$ws.Range("B6").Value = '[-, -, -, -]'
$ws.Range("C6").Value = '-'
$ws.Range("D6").Value = '-'
$ws.Range("E6").Value = '[-, ''MCT-3A-Processos de Usinagem 2'', -, -]'
$ws.Range("F6").Value = '[-, ''MEC-3B-Retífica'', -, -]'
$ws.Range("C7").Value = '-'
$ws.Range("D7").Value = '-'
$ws.Range("E7").Value = '[-, ''MCT-3A-Processos de Usinagem 2'', -, -]'
$ws.Range("F7").Value = '[-, ''MEC-3B-Retífica'', -, -]'
$ws.Range("C8").Value = '[-, -, -, -]'
$ws.Range("D8").Value = '[-, -, -, -]'
$ws.Range("B10").Value = 'MEC-3A-Elem. Máquinas'
$ws.Range("D10").Value = '[-, -, -, -]'
$ws.Range("E10").Value = '-'
$ws.Range("F14").Value = '[-, -, -, -]'
$ws.Range("B16").Value = 'MEC-3A-Elem. Máquinas'
$ws.Range("E16").Value = '-'
$ws.Range("B18").Value = '-'
$ws.Range("D18").Value = '-'
$ws.Range("E18").Value = '[-, ''ELM-1NA-Processos de Usinagem 2'', -, -]'
$ws.Range("F18").Value = '[''MEC-2NB-Retífica'', -, -, -]'
$ws.Range("B19").Value = '-'
$ws.Range("D19").Value = '-'
$ws.Range("E19").Value = '[-, ''ELM-1NA-Processos de Usinagem 2'', -, -]'
$ws.Range("F19").Value = '[''MEC-2NB-Retífica'', -, -, -]'
$ws.Range("B20").Value = '-'
$ws.Range("E20").Value = '[-, ''ELM-1NA-Processos de Usinagem 2'', -, -]'
$ws.Range("F20").Value = '[Joel L.-Fundição-2NB, -, -, -]'
$ws.Range("B21").Value = '-'
$ws.Range("E21").Value = '[-, ''ELM-1NA-Processos de Usinagem 2'', -, -]'
$ws.Range("F21").Value = '[''MEC-2NB-Retífica'', -, -, -]'
